{"js": "// Append three new dated journal-entry paragraphs to the end of the\n// document body (after the \"04-04-2023 ... roll over and come back up.\"\n// entry), matching the new content added by the commit:\n//   - a new \"04-19-2023\" date header\n//   - a status-update paragraph about the last few weeks of attempts\n//   - a paragraph about switching to PPO instead of DDPG\n\nconst body = context.document.body;\n\n// Insert each paragraph at the very end of the body (Word automatically\n// places new paragraphs before the body's trailing sectPr).\nconst p1 = body.insertParagraph(\n  \"04-19-2023 ---------------------------------------------------\",\n  Word.InsertLocation.end\n);\n\nconst p2 = body.insertParagraph(\n  \"Lots of attempts past few weeks. Not much improvement. Main change is doubled rotation torque to see if dynamics of system were limiting performance. This seems somewhat true since the latest attempts seem to achieve stability even when initialized upside down. The tracking performance seems worse though. \",\n  Word.InsertLocation.end\n);\n\nconst p3 = body.insertParagraph(\n  \"Next attempt is using Proximal Policy Optimization as the training agent rather than Deep Deterministic Policy Gradient. Some literature suggests better convergence which would be good since the latest sessions have been taking >24 hours \",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Append three new dated journal-entry paragraphs to the end of the\n# document body (after the \"04-04-2023 ... roll over and come back up.\"\n# entry), matching the new content added by the commit:\n#   - a new \"04-19-2023\" date header\n#   - a status-update paragraph about the last few weeks of attempts\n#   - a paragraph about switching to PPO instead of DDPG\n\n$d = $word.ActiveDocument\n\n$end = $d.Range()\n$end.Collapse(0)  # wdCollapseEnd\n\n$p1 = \"04-19-2023 ---------------------------------------------------\"\n$p2 = \"Lots of attempts past few weeks. Not much improvement. Main change is doubled rotation torque to see if dynamics of system were limiting performance. This seems somewhat true since the latest attempts seem to achieve stability even when initialized upside down. The tracking performance seems worse though. \"\n$p3 = \"Next attempt is using Proximal Policy Optimization as the training agent rather than Deep Deterministic Policy Gradient. Some literature suggests better convergence which would be good since the latest sessions have been taking >24 hours \"\n\n$end.InsertAfter(\"`r\" + $p1 + \"`r\" + $p2 + \"`r\" + $p3)\n"}
